# Common: Some another missing pieces added
#
# The "vendors" sheet (column A, sorted alphabetically/case-insensitively)
# is missing a vendor. Insert "Sirius Mods" in its correct sorted position
# (it sorts right after "Samsung" and before "SMArt Mods"), pushing the
# rest of the list down by one row, then re-assert the sheet's sort state
# over the full (now one-row-larger) range, and finally move the active
# selection to reflect where the user ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vendors")

# Insert a new row at 63 (between "Samsung" and "SMArt Mods") and fill it.
$ws.Rows("63:63").Insert()
$ws.Cells.Item(63, 1).Value = "Sirius Mods"

# Refresh the sort state/condition to cover the full A2:A94 list now that
# a row was added (mirrors re-running the existing column sort).
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("A2:A94"))
$sort.SetRange($ws.Range("A2:A94"))
$sort.Header = 0
$sort.Apply()

# Reflect the resulting selection on the sheet.
$ws.Range("A68").Select()
